# Update "想去人数" (interest count) values in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.
# Both sheets list the same events (全部类型 is a merged view of all
# category sheets), so the same F-column updates must be applied to the
# matching rows in each sheet.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value mapping for the "展览" sheet
$exhibitionUpdates = @{
    4  = 93
    8  = 26
    9  = 8251
    10 = 771
    11 = 272
    12 = 1115
    13 = 845
    14 = 49
    15 = 36
    16 = 209
    17 = 107
    20 = 891
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value mapping for the "全部类型" sheet
$allTypesUpdates = @{
    4  = 93
    10 = 26
    11 = 8251
    12 = 771
    13 = 272
    14 = 1115
    15 = 845
    16 = 49
    17 = 36
    18 = 209
    19 = 107
    22 = 891
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
